# Refresh cryptocurrency price / 1h-volume-change snapshot data (Price ==
# column D, Volume(1h) == column E) for the rows whose source values
# changed in this run of the GitHub Actions symbol-list scraper.
#
# Both columns hold plain text (e.g. "305.77", "5.54%"), not numbers --
# the leading "'" forces Excel to store the new value as text too,
# instead of silently reinterpreting it as a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'305.77"
$ws.Range("E2").Value = "'5.54%"

# Row 3: OKB
$ws.Range("D3").Value = "'32.31"
$ws.Range("E3").Value = "'9.48%"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'5.336"
$ws.Range("E4").Value = "'4.41%"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.07396"
$ws.Range("E5").Value = "'10.71%"

# Row 6: KuCoinToken
$ws.Range("D6").Value = "'7.757"
$ws.Range("E6").Value = "'5.22%"

# Row 7: GateToken
$ws.Range("D7").Value = "'3.700"
$ws.Range("E7").Value = "'8.53%"

# Row 8: FTXToken
$ws.Range("D8").Value = "'1.550"
$ws.Range("E8").Value = "'13.51%"

# Row 9: MXToken
$ws.Range("D9").Value = "'0.9101"
$ws.Range("E9").Value = "'-1.30%"

# Row 10: One
$ws.Range("D10").Value = "'0.01643"
$ws.Range("E10").Value = "'2,448.41%"

# Row 11: WazirX
$ws.Range("D11").Value = "'0.1669"
$ws.Range("E11").Value = "'4.98%"

# Row 12: LiechtensteinCryptoassetsExchange
$ws.Range("D12").Value = "'0.07530"
$ws.Range("E12").Value = "'12.31%"

# Row 13: MandalaExchangeToken
$ws.Range("D13").Value = "'0.08004"
$ws.Range("E13").Value = "'3.91%"

# Row 14: BitrueCoin
$ws.Range("D14").Value = "'0.03027"
$ws.Range("E14").Value = "'2.63%"

# Row 15: BitMartToken
$ws.Range("D15").Value = "'0.09863"
$ws.Range("E15").Value = "'9.66%"

# Row 16: BitForexToken
$ws.Range("D16").Value = "'0.001520"
$ws.Range("E16").Value = "'-3.31%"

# Row 17: CoinExToken
$ws.Range("D17").Value = "'0.04563"
$ws.Range("E17").Value = "'1.11%"

# Row 18: TigerCash
$ws.Range("D18").Value = "'0.006362"
$ws.Range("E18").Value = "'1.36%"

# Row 19: LEO
$ws.Range("D19").Value = "'3.475"
$ws.Range("E19").Value = "'0.71%"

# Row 20: BTSEToken
$ws.Range("D20").Value = "'2.239"
$ws.Range("E20").Value = "'0.87%"

# Row 21: BitpandaEcosystemToken
$ws.Range("E21").Value = "'1.52%"

# Row 22: ProBitToken
$ws.Range("E22").Value = "'1.47%"

# Row 23: MCDex
$ws.Range("D23").Value = "'4.218"
$ws.Range("E23").Value = "'3.33%"

# Row 24: ZBToken
$ws.Range("D24").Value = "'0.1631"
$ws.Range("E24").Value = "'4.06%"

# Row 25: BitKan
$ws.Range("D25").Value = "'0.001217"
$ws.Range("E25").Value = "'2.23%"

# Row 26: HotbitToken
$ws.Range("D26").Value = "'0.004502"
$ws.Range("E26").Value = "'8.86%"

# Row 27: NitroEx
$ws.Range("E27").Value = "'-6.36%"

# Row 28: UpBots
$ws.Range("D28").Value = "'0.0001804"
$ws.Range("E28").Value = "'11.48%"

# Row 40: IDEX
$ws.Range("E40").Value = "'6.51%"

# Row 41: KickToken
$ws.Range("D41").Value = "'0.007432"
$ws.Range("E41").Value = "'10.43%"

# Row 42: BKEXToken
$ws.Range("D42").Value = "'0.1362"
$ws.Range("E42").Value = "'9.78%"

# Row 43: CEJI
$ws.Range("E43").Value = "'14.19%"

# Row 44: LocalTraders
$ws.Range("D44").Value = "'0.01382"
$ws.Range("E44").Value = "'14.37%"

# Row 45: CoinLion
$ws.Range("D45").Value = "'0.00006153"
$ws.Range("E45").Value = "'7.94%"

# Row 46: BOLO
$ws.Range("D46").Value = "'1.893"
$ws.Range("E46").Value = "'-3.94%"

# Row 47: CoinbaseStockToken
$ws.Range("D47").Value = "'0.01302"
$ws.Range("E47").Value = "'-0.37%"
